$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two polyswitch/fuse rows (RUEF900-ND / F2 and RGEF500-ND / F1,F3) are no
# longer on the board ("no more airwires"), so remove rows 54:55 entirely.
# Everything below shifts up by two rows and the SUM() formulas get
# automatically re-targeted by Excel.
$ws.Range("A54:J55").EntireRow.Delete()

$ws.Range("C54").Select()
